$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws2 = $wb.Worksheets.Item("演出")
$ws3 = $wb.Worksheets.Item("本地生活")
$ws4 = $wb.Worksheets.Item("全部类型")

$ws1.Range("F5").Value = 1241  # was 1238
$ws1.Range("F7").Value = 7489  # was 7482
$ws1.Range("F11").Value = 8157  # was 8150
$ws1.Range("F13").Value = 58  # was 57
$ws1.Range("F14").Value = 5573  # was 5566
$ws1.Range("F16").Value = 2530  # was 2527
$ws1.Range("F17").Value = 1082  # was 1079
$ws1.Range("F18").Value = 4574  # was 4572
$ws1.Range("F19").Value = 317  # was 316
$ws1.Range("F21").Value = 88  # was 87
$ws1.Range("F23").Value = 450  # was 448
$ws1.Range("F24").Value = 1672  # was 1619
$ws1.Range("F25").Value = 28  # was 27
$ws1.Range("F26").Value = 2667  # was 2640
$ws1.Range("F28").Value = 306  # was 303
$ws1.Range("F29").Value = 105  # was 104
$ws1.Range("F30").Value = 240  # was 235
$ws1.Range("F31").Value = 619  # was 618
$ws1.Range("F33").Value = 313  # was 310
$ws1.Range("F34").Value = 1598  # was 1592
$ws1.Range("F37").Value = 2538  # was 2530
$ws1.Range("F38").Value = 2257  # was 2252
$ws1.Range("F39").Value = 6  # was 4
$ws2.Range("F3").Value = 99  # was 96
$ws2.Range("F4").Value = 29  # was 28
$ws3.Range("F3").Value = 1288  # was 1286
$ws4.Range("F4").Value = 1288  # was 1286
$ws4.Range("F7").Value = 1241  # was 1238
$ws4.Range("F9").Value = 7489  # was 7482
$ws4.Range("F13").Value = 8157  # was 8151
$ws4.Range("F15").Value = 58  # was 57
$ws4.Range("F16").Value = 5573  # was 5566
$ws4.Range("F18").Value = 2530  # was 2527
$ws4.Range("F19").Value = 1082  # was 1079
$ws4.Range("F20").Value = 4574  # was 4572
$ws4.Range("F21").Value = 317  # was 316
$ws4.Range("F23").Value = 88  # was 87
$ws4.Range("F26").Value = 99  # was 96
$ws4.Range("F27").Value = 450  # was 448
$ws4.Range("F28").Value = 1672  # was 1620
$ws4.Range("F29").Value = 28  # was 27
$ws4.Range("F30").Value = 2667  # was 2640
$ws4.Range("F32").Value = 306  # was 303
$ws4.Range("F33").Value = 105  # was 104
$ws4.Range("F34").Value = 240  # was 235
$ws4.Range("F35").Value = 29  # was 28
$ws4.Range("F36").Value = 619  # was 618
$ws4.Range("F38").Value = 313  # was 310
$ws4.Range("F40").Value = 1598  # was 1592
$ws4.Range("F43").Value = 2538  # was 2530
$ws4.Range("F45").Value = 2257  # was 2252
$ws4.Range("F46").Value = 6  # was 4
